$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Row 6 on ProductLoanInput: relabel the currency field and drop the
# trailing space from its value, restyling the value cell to match the
# other "accounting/cash" highlighted input cells (green fill, normal font).
$ws1.Cells.Item(6, 1).Value2 = "currency"

$valueCell = $ws1.Cells.Item(6, 2)
$valueCell.Style = "Normal"
$valueCell.Value2 = "US Dollar"
$valueCell.Interior.Color = 5296274

# Make ProductLoanInput the active sheet/tab with A6:B6 selected.
$ws1.Select()
$ws1.Range("A6:B6").Select()
